# Apply "nota EV02 y sprint 2" updates to Sprint_Backlog_Actividades_Horas_DFF.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product_Backlog_tareas")

# Update estimated-hours (column D) values for the affected tasks.
$ws.Range("D9").Value  = 24
$ws.Range("D11").Value = 10
$ws.Range("D12").Value = 8
$ws.Range("D13").Value = 8
$ws.Range("D14").Value = 7
$ws.Range("D15").Value = 9
$ws.Range("D17").Value = 29
$ws.Range("D18").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("D21").Value = 10
$ws.Range("D24").Value = 28
$ws.Range("D25").Value = 2

# Recalculate so the Total Horas formula (C2 = SUM(D5:D78)) reflects the new values.
$excel.Calculate()

# Scroll the sheet view to match where the author was last working and
# update the active selection, as captured in the saved view state.
$window = $excel.ActiveWindow
$window.ScrollRow = 30
$window.ScrollColumn = 1
$ws.Range("D34").Select()

$wb.Save()
